# Fix imported idProveedor value on row 2 (was mistakenly -1, should be 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 10).Value = 1

# Update the saved view state: zoom to 53% and move the selection to E9
# (this also clears the old topLeftCell="F1" scroll anchor)
$excel.ActiveWindow.Zoom = 53
$ws.Range("E9").Select()

# Columns J (idProveedor) and N (observaciones) now carry an explicit
# best-fit style width, matching the widths already used by columns
# M (fotoProducto) and I (fechaIngreso) respectively.
$ws.Columns.Item(10).ColumnWidth = 13.85546875
$ws.Columns.Item(14).ColumnWidth = 15.7109375
